$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3707.1538
$ws.Range("J76").Value = 3739.1
$ws.Range("L76").Value = 3739.1
$ws.Range("N76").Value = -4369.1
$ws.Range("H79").Value = 3707.1538
$ws.Range("J79").Value = 3739.1
$ws.Range("L79").Value = 3739.1
$ws.Range("N79").Value = -5923.1
$ws.Range("H106").Value = 2643.697
$ws.Range("I106").Value = 4238.5
$ws.Range("J106").Value = 1732.381
$ws.Range("K106").Value = 4238.5
$ws.Range("L106").Value = 1732.381
$ws.Range("M106").Value = -3607.5
$ws.Range("N106").Value = -2994.381
$ws.Range("H107").Value = 2139.1428
$ws.Range("I107").Value = 2380.4
$ws.Range("K107").Value = 2380.4
$ws.Range("M107").Value = -460.4000000000001
$ws.Range("H132").Value = 6155.683
$ws.Range("I132").Value = 6388.7295
$ws.Range("K132").Value = 19166.1885
$ws.Range("M132").Value = -16636.1885
$ws.Range("H135").Value = 5897.7144
$ws.Range("I135").Value = 1321
$ws.Range("K135").Value = 11889
$ws.Range("M135").Value = -9354
$ws.Range("H137").Value = 21747598
$ws.Range("I137").Value = 31252386
$ws.Range("J137").Value = 22367.857
$ws.Range("K137").Value = 93757158
$ws.Range("L137").Value = 67103.571
$ws.Range("M137").Value = -93754608
$ws.Range("N137").Value = -72203.571
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 977608.25
$ws.Range("I32").Value = 1090368.1
$ws.Range("J32").Value = 19149.834
$ws.Range("K32").Value = 1090368.1
$ws.Range("L32").Value = 19149.834
$ws.Range("M32").Value = -1090081.1
$ws.Range("N32").Value = -19723.834
$ws.Range("H45").Value = 1725.826
$ws.Range("I45").Value = 1555.6471
$ws.Range("J45").Value = 2208
$ws.Range("K45").Value = 1555.6471
$ws.Range("L45").Value = 2208
$ws.Range("M45").Value = -1178.6471
$ws.Range("N45").Value = -2962
$ws.Range("H61").Value = 4008482.5
$ws.Range("I61").Value = 10265.277
$ws.Range("J61").Value = 14289613
$ws.Range("K61").Value = 10265.277
$ws.Range("L61").Value = 14289613
$ws.Range("M61").Value = -10053.277
$ws.Range("N61").Value = -14290037
$ws.Range("H110").Value = 2253.3845
$ws.Range("I110").Value = 3098
$ws.Range("K110").Value = 3098
$ws.Range("M110").Value = -1053
$ws.Range("H132").Value = 5982.32
$ws.Range("I132").Value = 3851.2727
$ws.Range("K132").Value = 11553.8181
$ws.Range("M132").Value = -9023.8181
$ws.Range("H136").Value = 4008482.5
$ws.Range("I136").Value = 10265.277
$ws.Range("J136").Value = 14289613
$ws.Range("K136").Value = 30795.831
$ws.Range("L136").Value = 42868839
$ws.Range("M136").Value = -28245.831
$ws.Range("N136").Value = -42873939
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2593.4375
$ws.Range("I86").Value = 2408.75
$ws.Range("J86").Value = 2901.25
$ws.Range("K86").Value = 2408.75
$ws.Range("L86").Value = 2901.25
$ws.Range("M86").Value = -1285.75
$ws.Range("N86").Value = -5147.25
$ws.Range("H89").Value = 2593.4375
$ws.Range("I89").Value = 2408.75
$ws.Range("J89").Value = 2901.25
$ws.Range("K89").Value = 12043.75
$ws.Range("L89").Value = 14506.25
$ws.Range("M89").Value = -6427.75
$ws.Range("N89").Value = -25738.25
$ws.Range("H94").Value = 2796.4583
$ws.Range("I94").Value = 2796.4583
$ws.Range("K94").Value = 2796.4583
$ws.Range("M94").Value = -2345.4583
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 971.06665
$ws.Range("I107").Value = 731.0909
$ws.Range("K107").Value = 731.0909
$ws.Range("M107").Value = 1188.9091
$ws.Range("H132").Value = 6293.3335
$ws.Range("I132").Value = 6642.5
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 19927.5
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -17397.5
$ws.Range("N132").Value = -15560
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2999.5
$ws.Range("J52").Value = 2999.5
$ws.Range("L52").Value = 8998.5
$ws.Range("N52").Value = -9530.5
$ws.Range("H92").Value = 215
$ws.Range("I92").Value = 213.33333
$ws.Range("J92").Value = 216.66667
$ws.Range("K92").Value = 639.99999
$ws.Range("L92").Value = 650.00001
$ws.Range("M92").Value = 608.00001
$ws.Range("N92").Value = -3146.00001
$ws.Range("H109").Value = 5453.4
$ws.Range("I109").Value = 1452.8572
$ws.Range("K109").Value = 4358.571599999999
$ws.Range("M109").Value = -3318.571599999999
$ws.Range("H124").Value = 11347.818
$ws.Range("I124").Value = 8304.333000000001
$ws.Range("K124").Value = 24912.999
$ws.Range("M124").Value = -20002.999
$ws.Range("H129").Value = 2131.875
$ws.Range("I129").Value = 644.1429000000001
$ws.Range("J129").Value = 3289
$ws.Range("K129").Value = 1932.4287
$ws.Range("L129").Value = 9867
$ws.Range("M129").Value = 3067.5713
$ws.Range("N129").Value = -19867
$ws.Range("H131").Value = 7660.8335
$ws.Range("I131").Value = 1763.1111
$ws.Range("J131").Value = 11199.467
$ws.Range("K131").Value = 5289.3333
$ws.Range("L131").Value = 33598.401
$ws.Range("M131").Value = -249.3333000000002
$ws.Range("N131").Value = -43678.401
$ws.Range("H132").Value = 1865.3125
$ws.Range("I132").Value = 1694.5454
$ws.Range("K132").Value = 15250.9086
$ws.Range("M132").Value = -12720.9086
$ws.Range("H137").Value = 7496.7856
$ws.Range("I137").Value = 2814.8333
$ws.Range("J137").Value = 11008.25
$ws.Range("K137").Value = 8444.499899999999
$ws.Range("L137").Value = 33024.75
$ws.Range("M137").Value = -3344.499899999999
$ws.Range("N137").Value = -43224.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1412.3334
$ws.Range("I80").Value = 1325
$ws.Range("K80").Value = 1325
$ws.Range("M80").Value = -327
$ws.Range("H83").Value = 1412.3334
$ws.Range("I83").Value = 1325
$ws.Range("K83").Value = 6625
$ws.Range("M83").Value = -1633
$ws.Range("H102").Value = 1981.4286
$ws.Range("I102").Value = 981.1667
$ws.Range("K102").Value = 981.1667
$ws.Range("M102").Value = 640.8333
$ws.Range("H113").Value = 1558.3334
$ws.Range("I113").Value = 1590.909
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 1590.909
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 579.0909999999999
$ws.Range("N113").Value = -5540
$ws.Range("H122").Value = 39365.15
$ws.Range("I122").Value = 51892.55
$ws.Range("K122").Value = 155677.65
$ws.Range("M122").Value = -153227.65
$ws.Range("H132").Value = 9643.423000000001
$ws.Range("I132").Value = 6382.5
$ws.Range("K132").Value = 19147.5
$ws.Range("M132").Value = -16617.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 13333.333
$ws.Range("J59").Value = 13333.333
$ws.Range("L59").Value = 13333.333
$ws.Range("N59").Value = -14641.333
$ws.Range("H82").Value = 97186.75
$ws.Range("I82").Value = 69330.664
$ws.Range("K82").Value = 69330.664
$ws.Range("M82").Value = -68969.664
$ws.Range("H85").Value = 97186.75
$ws.Range("I85").Value = 69330.664
$ws.Range("K85").Value = 69330.664
$ws.Range("M85").Value = -68082.664
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1298.5264
$ws.Range("I107").Value = 1078.2667
$ws.Range("K107").Value = 3234.800099999999
$ws.Range("M107").Value = -1314.800099999999
$ws.Range("H122").Value = 40207.965
$ws.Range("I122").Value = 1143.0952
$ws.Range("K122").Value = 3429.2856
$ws.Range("M122").Value = -979.2856000000002
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("N131").Value = 0
$ws.Range("H132").Value = 5557354
$ws.Range("I132").Value = 6668254
$ws.Range("K132").Value = 20004762
$ws.Range("M132").Value = -20002232
